$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J, row 1 (matches style of existing header cells)
$ws.Range("J1").Value = "E-Mail Dozentenumfrage Aufwandbereitschaft"

# Fill column J rows 2-6 with "X" markers, matching the other criteria columns
$ws.Range("J2").Value = "X"
$ws.Range("J3").Value = "X"
$ws.Range("J4").Value = "X"
$ws.Range("J5").Value = "X"
$ws.Range("J6").Value = "X"

# New "Ausgefüllt von" entry for column J, row 8 (match formatting of I8)
$ws.Range("I8").Copy() | Out-Null
$ws.Range("J8").PasteSpecial(-4122) | Out-Null
$ws.Range("J8").Value = "Heiser/ Netzler"

# Update the active selection to K8, reflecting the new end-of-entry cell
$ws.Range("K8").Select()
